# Applies the "Lecture partielle de l'EDT M1 MIAGE" edit:
#  - The day-of-week label "lundi" becomes "vendredi"
#  - The five week-start dates (col A) are shifted from Monday dates
#    in 2023 to the corresponding Friday dates in 2026

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

# Update the weekday label used in column B (shared string "lundi" -> "vendredi")
for ($r = 2; $r -le 14; $r += 3) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "lundi") {
        $cell.Value2 = "vendredi"
    }
}

# Update the dates in column A (stored as date-formatted serial numbers,
# keep the time-of-day at midnight so the serial stays a whole number)
$ws.Cells.Item(2, 1).Value2  = Get-Date -Year 2026 -Month 2 -Day 6  -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(5, 1).Value2  = Get-Date -Year 2026 -Month 2 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(8, 1).Value2  = Get-Date -Year 2026 -Month 2 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(11, 1).Value2 = Get-Date -Year 2026 -Month 2 -Day 27 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(14, 1).Value2 = Get-Date -Year 2026 -Month 3 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
